$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D12: clear the title text (make it empty)
$ws.Range("D12").ClearContents()

# E12: update link
$ws.Range("E12").Value = "https://tensorflow.blog/2023/10/20/book-roadmap/"

# D51: update title
$ws.Range("D51").Value = "[aws] RDS 타임존 변경하는 방법"

# E51: update link
$ws.Range("E51").Value = "https://bskyvision.com/entry/aws-RDS-%ED%83%80%EC%9E%84%EC%A1%B4-%EB%B3%80%EA%B2%BD%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
